$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.112942218780518
$ws.Range("B1").Value = 4.194282054901123
$ws.Range("C1").Value = 2.106475591659546
$ws.Range("D1").Value = 1.625145673751831
$ws.Range("E1").Value = 1.460566997528076
